$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Hate Crimes row: C31 changes from numeric 1 to text "0" (same style as other N/A cells) ---
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Cells.Item(22, 14).Copy()
$ws.Range("C31").PasteSpecial(-4122)

# --- Numeric data updates across rows 14-33 ---
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 600
$ws.Range("I14").Value = 63
$ws.Range("K14").Value = 10.526315789473
$ws.Range("L14").Value = -10
$ws.Range("M14").Value = -47.5
$ws.Range("N14").Value = -84.855769230769
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -33.333333333333
$ws.Range("G15").Value = 17
$ws.Range("H15").Value = -17.647058823529
$ws.Range("I15").Value = 214
$ws.Range("J15").Value = 192
$ws.Range("K15").Value = 11.458333333333
$ws.Range("L15").Value = -5.309734513274
$ws.Range("M15").Value = 5.940594059405
$ws.Range("N15").Value = -59.546313799621
$ws.Range("C16").Value = 46
$ws.Range("D16").Value = 47
$ws.Range("E16").Value = -2.127659574468
$ws.Range("F16").Value = 195
$ws.Range("G16").Value = 239
$ws.Range("H16").Value = -18.410041841004
$ws.Range("I16").Value = 2075
$ws.Range("J16").Value = 2169
$ws.Range("K16").Value = -4.333794375288
$ws.Range("L16").Value = -7.200357781753
$ws.Range("M16").Value = -34.459886291850
$ws.Range("N16").Value = -85.660977126667
$ws.Range("C17").Value = 67
$ws.Range("D17").Value = 71
$ws.Range("E17").Value = -5.633802816901
$ws.Range("G17").Value = 304
$ws.Range("H17").Value = -3.618421052631
$ws.Range("I17").Value = 3754
$ws.Range("J17").Value = 3703
$ws.Range("K17").Value = 1.377261679719
$ws.Range("L17").Value = 3.160208848584
$ws.Range("M17").Value = 29.941156109380
$ws.Range("N17").Value = -49.637778374027
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = 71.428571428571
$ws.Range("F18").Value = 130
$ws.Range("G18").Value = 133
$ws.Range("H18").Value = -2.255639097744
$ws.Range("I18").Value = 1647
$ws.Range("J18").Value = 1763
$ws.Range("K18").Value = -6.579693703913
$ws.Range("L18").Value = -21.534063839923
$ws.Range("M18").Value = -40.670028818443
$ws.Range("N18").Value = -84.251290877796
$ws.Range("C19").Value = 114
$ws.Range("D19").Value = 121
$ws.Range("E19").Value = -5.785123966942
$ws.Range("F19").Value = 422
$ws.Range("G19").Value = 470
$ws.Range("H19").Value = -10.212765957446
$ws.Range("I19").Value = 4641
$ws.Range("J19").Value = 5072
$ws.Range("K19").Value = -8.497634069400
$ws.Range("L19").Value = -11.278914165551
$ws.Range("M19").Value = 24.959612277867
$ws.Range("N19").Value = -21.802864363942
$ws.Range("C20").Value = 36
$ws.Range("D20").Value = 38
$ws.Range("E20").Value = -5.263157894736
$ws.Range("F20").Value = 116
$ws.Range("H20").Value = -21.621621621621
$ws.Range("I20").Value = 1520
$ws.Range("J20").Value = 1593
$ws.Range("K20").Value = -4.582548650345
$ws.Range("L20").Value = -5.059337913803
$ws.Range("M20").Value = 22.481869460112
$ws.Range("N20").Value = -81.485992691839
$ws.Range("D21").Value = 301
$ws.Range("E21").Value = 0.664451827242
$ws.Range("F21").Value = 1177
$ws.Range("G21").Value = 1312
$ws.Range("H21").Value = -10.289634146341
$ws.Range("I21").Value = 13914
$ws.Range("J21").Value = 14549
$ws.Range("K21").Value = -4.364561138222
$ws.Range("L21").Value = -7.866507747318
$ws.Range("M21").Value = -1.375106322653
$ws.Range("N21").Value = -70.690708402671
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -11.111111111111
$ws.Range("F22").Value = 22
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = -21.428571428571
$ws.Range("I22").Value = 251
$ws.Range("J22").Value = 258
$ws.Range("K22").Value = -2.713178294573
$ws.Range("L22").Value = -17.973856209150
$ws.Range("M22").Value = -33.947368421052
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 21
$ws.Range("E23").Value = -9.523809523809
$ws.Range("F23").Value = 95
$ws.Range("G23").Value = 102
$ws.Range("H23").Value = -6.862745098039
$ws.Range("I23").Value = 1238
$ws.Range("J23").Value = 1361
$ws.Range("K23").Value = -9.037472446730
$ws.Range("L23").Value = -6.283118849356
$ws.Range("M23").Value = 21.135029354207
$ws.Range("C24").Value = 273
$ws.Range("D24").Value = 178
$ws.Range("E24").Value = 53.370786516853
$ws.Range("F24").Value = 1078
$ws.Range("G24").Value = 822
$ws.Range("H24").Value = 31.143552311435
$ws.Range("I24").Value = 10824
$ws.Range("J24").Value = 10597
$ws.Range("K24").Value = 2.142115693120
$ws.Range("L24").Value = -7.336700624946
$ws.Range("M24").Value = 19.681556833259
$ws.Range("C25").Value = 126
$ws.Range("D25").Value = 69
$ws.Range("E25").Value = 82.608695652173
$ws.Range("F25").Value = 463
$ws.Range("G25").Value = 334
$ws.Range("H25").Value = 38.622754491018
$ws.Range("I25").Value = 4822
$ws.Range("J25").Value = 4219
$ws.Range("K25").Value = 14.292486371178
$ws.Range("L25").Value = -1.672104404567
$ws.Range("C26").Value = 129
$ws.Range("D26").Value = 106
$ws.Range("E26").Value = 21.698113207547
$ws.Range("F26").Value = 516
$ws.Range("H26").Value = 5.737704918032
$ws.Range("I26").Value = 5562
$ws.Range("J26").Value = 5307
$ws.Range("K26").Value = 4.804974561899
$ws.Range("L26").Value = 8.210116731517
$ws.Range("M26").Value = -19.53125
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = -60
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 25
$ws.Range("H27").Value = -28
$ws.Range("I27").Value = 298
$ws.Range("J27").Value = 295
$ws.Range("K27").Value = 1.016949152542
$ws.Range("L27").Value = -11.309523809523
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 49
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = 19.512195121951
$ws.Range("I28").Value = 575
$ws.Range("J28").Value = 546
$ws.Range("K28").Value = 5.311355311355
$ws.Range("L28").Value = 6.877323420074
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 18
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = 20
$ws.Range("I29").Value = 213
$ws.Range("J29").Value = 206
$ws.Range("K29").Value = 3.398058252427
$ws.Range("L29").Value = -30.844155844155
$ws.Range("M29").Value = -53.493449781659
$ws.Range("N29").Value = -87.191822008418
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 17
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = 21.428571428571
$ws.Range("I30").Value = 176
$ws.Range("J30").Value = 175
$ws.Range("K30").Value = 0.571428571428
$ws.Range("L30").Value = -31.782945736434
$ws.Range("M30").Value = -52.688172043010
$ws.Range("N30").Value = -88.219544846050
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 15
$ws.Range("G31").Value = 17
$ws.Range("H31").Value = -11.764705882352
$ws.Range("I31").Value = 79
$ws.Range("J31").Value = 71
$ws.Range("K31").Value = 11.267605633802
$ws.Range("L31").Value = -2.469135802469
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -25
$ws.Range("I33").Value = 24
$ws.Range("J33").Value = 20
$ws.Range("K33").Value = 20
$ws.Range("L33").Value = -7.692307692307
